$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the cable length and quantity values (e.g. new CBL14_24V_open variant halves length/qty)
$ws.Range("C15").Value = 50
$ws.Range("C22").Value = 50

# Move the active selection as reflected by the saved view state
$ws.Range("E14").Select()
